# Updates cryptocurrency price (D) and 1h volume-change (E) columns.
# Values are assigned with a leading apostrophe and the style is reset to
# "Normal" afterward so numeric-looking strings (e.g. "0.999", "331.05")
# are stored as text (matching column D/E's existing inline-string cells)
# instead of being auto-coerced to numbers, without leaving a NumberFormat
# / style footprint on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.390.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.69%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.741.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.90%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.12%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'115.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.32%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'331.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.29%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.42%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.12%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.564"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.20%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'41.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.31%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'20.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.10%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.0829"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.41%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +2.67%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +3.44%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.156.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.40%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.715.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.41%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.885"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.60%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'51.246.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.46%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'13.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.33%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'3.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +3.13%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.72%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.0₃0962"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.29%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'289.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +4.31%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'70.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.06%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.19%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'26.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.10%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.02%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'10.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.98%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.36%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.142"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.18%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'35.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.81%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'50.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.58%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'5.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.55%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +1.06%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'19.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.30%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.15%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'5.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.93%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.48%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.69%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'24.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +9.25%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'129.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.92%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.0352"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +10.24%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +3.55%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.02%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'3.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.95%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.117.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.30%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +11.58%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.21%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'5.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.24%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'9.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.55%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'60.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.89%  "
$ws.Range("E51").Style = "Normal"
